$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 2020 column (Q) values for rows 4 (year header) through 14 (Kyrgyz Republic total)
$values = @{
  4  = 2020
  5  = 109.7221295941265
  6  = 108.44905375816947
  7  = 109.90982951756889
  8  = 108.40606487500015
  9  = 109.40161876466024
  10 = 107.71155656686271
  11 = 111.78921596090774
  12 = 111.39254046803097
  13 = 110.44919152842827
  14 = 106.89826464456031
}

foreach ($row in $values.Keys) {
  $ws.Cells.Item($row, 17).Value = $values[$row]
}

# Copy the formatting of column P (the preceding year column) onto the
# newly populated column Q cells, row by row, so each new cell matches the
# style used by its neighbour.
for ($row = 4; $row -le 14; $row++) {
  $ws.Cells.Item($row, 16).Copy() | Out-Null
  $ws.Cells.Item($row, 17).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Restore the selection to the cell the author left active (N14) so the
# saved sheetView reflects that.
$ws.Range("N14").Select() | Out-Null
